$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update a few existing "User" cells ---
$ws.Range("B3").Value = "SWG"
$ws.Range("B4").Value = "SWG"
$ws.Range("B10").Value = "internal/Imaging"

# --- Insert two new rows before the old row 19 ---
$ws.Range("A19:A20").EntireRow.Insert()

# Copy formatting of the blank separator row (row 18) down into the two
# freshly inserted rows so they pick up the same body-row style used
# throughout the table.
$ws.Range("A18:H18").Copy()
$ws.Range("A19:H20").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Rows.Item(19).RowHeight = 15
$ws.Rows.Item(20).RowHeight = 15

# --- Populate the new rows ---
$ws.Range("A19").Value = 2
$ws.Range("B19").Value = "internal"
$ws.Range("E19").Value = "med"
$ws.Range("F19").Value = 10
$ws.Range("H19").Value = "Update Introduce authorization components (gridgrouper,  csm, etc) to new authorization support"

$ws.Range("A20").Value = 2
$ws.Range("B20").Value = "internal"
$ws.Range("E20").Value = "med"
$ws.Range("F20").Value = 20
$ws.Range("H20").Value = "Add authorization configuration support in Introduce to Resources (needed for stateful services, such as FQP, workflow, BDT, etc)"

# --- Page setup: fit to page, scale 70%, landscape ---
$ws.PageSetup.Zoom = 70
$ws.PageSetup.Orientation = 2
$ws.PageSetup.FitToPagesWide = 1
$ws.PageSetup.FitToPagesTall = 1

# --- View: select H21 (the row that used to be row 19) and scroll to top ---
$ws.Range("H21").Select()

Write-Host "done"
